# The document contains the literal text "<id>p072v_1</id>" split across
# three separate runs:
#   1. "<id>"     - Courier New, color 7f6000, sz/szCs 18
#   2. "p072v_1"  - default font, color 000000
#   3. "</id>"    - Courier New, color 7f6000, sz/szCs 18
#
# The edit collapses these three runs into a single run containing the
# whole string "<id>p072v_1</id>", taking on the formatting of the first
# run in the range (Courier New / 7f6000 / 18pt), which matches the first
# and third runs' shared formatting.
#
# Find & Replace across the range naturally merges the matched runs into
# one run carrying the formatting of the first run found, which is exactly
# the desired outcome here.

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p072v_1</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p072v_1</id>", 2)
